$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GPAC")

$ws.Range("B5").Value = 92
$ws.Range("B6").Value = 238.05078783332203
$ws.Range("B7").Value = 17.239999999999998
